$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 21; existing rows 21-105 shift down to 22-106
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with its values.
$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value = "Bíobío"
$ws.Cells.Item(21, 4).Value = 44972
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112012
$ws.Cells.Item(21, 7).Value = "Espinaca"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 50
$ws.Cells.Item(21, 11).Value = 7000
$ws.Cells.Item(21, 12).Value = 7500
$ws.Cells.Item(21, 13).Value = 7200
$ws.Cells.Item(21, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(21, 15).Value = "Región Metropolitana"
$ws.Cells.Item(21, 16).Value = 720
$ws.Cells.Item(21, 17).Value = 10
$ws.Cells.Item(21, 18).Value = "Hortaliza"
